$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the language column for sub-rows that belong to the same
# language group as the row above them.
$ws.Range("A7").Value = $ws.Range("A6").Text
$ws.Range("A8").Value = $ws.Range("A6").Text
$ws.Range("A19").Value = $ws.Range("A18").Text
$ws.Range("A21").Value = $ws.Range("A20").Text
$ws.Range("A25").Value = $ws.Range("A24").Text

# Update the active selection shown when the sheet was last saved.
$ws.Range("B14").Select()
